# Inserts a new data row at row 270 in the "Papa" price list sheet.
# Every existing row from 270..311 shifts down by one (to 271..312),
# and the new row 270 carries a new price observation (same market/
# category metadata as the row that used to sit at 270, now at 271,
# but with its own date, variety, quality, volume and origin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 270:311 down to 271:312, leaving a blank row 270.
$ws.Rows.Item(270).Insert()

# Populate the new row 270.
$ws.Cells.Item(270, 1).Value = 5
$ws.Cells.Item(270, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(270, 3).Value = "Maule"
$ws.Cells.Item(270, 4).Value = 44505
$ws.Cells.Item(270, 5).Value = 7
$ws.Cells.Item(270, 6).Value = 100114001
$ws.Cells.Item(270, 7).Value = "Papa"
$ws.Cells.Item(270, 8).Value = "Rodeo"
$ws.Cells.Item(270, 9).Value = "1a nueva(o)"
$ws.Cells.Item(270, 10).Value = 1600
$ws.Cells.Item(270, 11).Value = 10000
$ws.Cells.Item(270, 12).Value = 10000
$ws.Cells.Item(270, 13).Value = 10000
$ws.Cells.Item(270, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(270, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(270, 16).Value = 400
$ws.Cells.Item(270, 17).Value = 25
$ws.Cells.Item(270, 18).Value = "Hortaliza"
